$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Column D = Price, Column E = Volume(1h)
# Each entry: RowNumber, NewPriceValue (or $null if unchanged), NewVolumeValue (or $null if unchanged)
# Price values are plain digit/dot strings that Excel's COM layer would happily
# reinterpret as numbers (e.g. "213.03", "63.17"). The source workbook keeps
# these as text cells, so we force text entry with a leading apostrophe and
# then reset the cell style back to Normal (the apostrophe trick stamps a
# "quote prefix" style onto the cell, which we don't want to keep).
$updates = @(
    @(2,  "28.442.45", "  +3.40%  "),
    @(3,  "1.591.86",  "  +1.44%  "),
    @(5,  "213.03",    "  +0.80%  "),
    @(6,  "0.491",     "  +0.09%  "),
    @(7,  $null,       "  +0.97%  "),
    @(8,  "24.41",     "  +7.49%  "),
    @(9,  $null,       "  +0.33%  "),
    @(10, $null,       "  +0.83%  "),
    @(11, $null,       "  +1.70%  "),
    @(12, "1.818.17",  "  +1.46%  "),
    @(13, "1.589.69",  "  +1.79%  "),
    @(14, $null,       "  +2.13%  "),
    @(15, $null,       "  -0.22%  "),
    @(16, "28.451.37", "  +3.58%  "),
    @(17, $null,       "  +1.21%  "),
    @(18, "229.88",    "  +1.68%  "),
    @(19, "0.0₃0707",  "  +0.08%  "),
    @(20, $null,       "  -0.58%  "),
    @(21, $null,       "  +0.84%  "),
    @(22, $null,       "  -1.42%  "),
    @(23, $null,       "  -0.73%  "),
    @(25, "151.71",    $null),
    @(26, "15.24",     "  +0.49%  "),
    @(27, $null,       "  -0.73%  "),
    @(28, $null,       "  -0.61%  "),
    @(29, $null,       "  +0.96%  "),
    @(30, $null,       "  -0.76%  "),
    @(31, $null,       "  +0.36%  "),
    @(32, $null,       "  +0.31%  "),
    @(33, $null,       "  +0.91%  "),
    @(34, "1.403.33",  "  -3.49%  "),
    @(35, $null,       "  -0.90%  "),
    @(36, $null,       "  -9.79%  "),
    @(37, $null,       "  +0.97%  "),
    @(38, "2.59",      "  +9.14%  "),
    @(39, $null,       "  -0.59%  "),
    @(40, "0.541",     "  +0.37%  "),
    @(41, "0.813",     "  -0.04%  "),
    @(42, $null,       "  +0.94%  "),
    @(43, "5.62",      "  -2.26%  "),
    @(44, "1.87",      "  +0.57%  "),
    @(45, "0.983",     "  +0.85%  "),
    @(46, "63.17",     "  -1.82%  "),
    @(47, "1.727.52",  "  +1.37%  "),
    @(48, $null,       "  +1.76%  "),
    @(49, $null,       "  +0.36%  "),
    @(50, $null,       "  +1.75%  "),
    @(51, $null,       "  -0.87%  ")
)

foreach ($row in $updates) {
    $rowNum = $row[0]
    $priceVal = $row[1]
    $volVal = $row[2]

    if ($null -ne $priceVal) {
        $cell = $ws.Cells.Item($rowNum, 4)
        $cell.Value = "'" + $priceVal
        $cell.Style = "Normal"
    }
    if ($null -ne $volVal) {
        $ws.Cells.Item($rowNum, 5).Value = $volVal
    }
}
